$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update deaths (H), reward (I), real (J) columns for rows 2-7
$ws.Range("H2").Value = 35746630268.87081
$ws.Range("I2").Value = 10735.6693187202
$ws.Range("J2").Value = 32204521182.33292

$ws.Range("H3").Value = 28549297387.41831
$ws.Range("I3").Value = 9754.205102515734
$ws.Range("J3").Value = 25331010990.93549

$ws.Range("H4").Value = 51911235557.14983
$ws.Range("I4").Value = 72870.41979858115
$ws.Range("J4").Value = 27868489041.26198

$ws.Range("H5").Value = 35746630268.87081
$ws.Range("I5").Value = 10735.6693187202
$ws.Range("J5").Value = 21468851863.61272

$ws.Range("H6").Value = 28549297387.41831
$ws.Range("I6").Value = 9754.205102515734
$ws.Range("J6").Value = 15576805888.41975

$ws.Range("H7").Value = 51911235557.14983
$ws.Range("I7").Value = 72870.41979858115
$ws.Range("J7").Value = -45001930757.31915
